$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 68.416664
$ws.Range("I5").Value = 61.22222
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 61.22222
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 53.77778
$ws.Range("N5").Value = -320

$ws.Range("H12").Value = 7355.6665
$ws.Range("I12").Value = 2600.5
$ws.Range("J12").Value = 8714.286
$ws.Range("K12").Value = 2600.5
$ws.Range("L12").Value = 8714.286
$ws.Range("M12").Value = -2430.5
$ws.Range("N12").Value = -9054.286

$ws.Range("H19").Value = 296.4
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 296.4
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 296.4
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -646.4

$ws.Range("H32").Value = 896.1667
$ws.Range("I32").Value = 675
$ws.Range("J32").Value = 1006.75
$ws.Range("K32").Value = 675
$ws.Range("L32").Value = 1006.75
$ws.Range("M32").Value = -349
$ws.Range("N32").Value = -1658.75

$ws.Range("H111").Value = 1439.125
$ws.Range("I111").Value = 1343.7142
$ws.Range("J111").Value = 1513.3334
$ws.Range("K111").Value = 4031.1426
$ws.Range("L111").Value = 4540.0002
$ws.Range("M111").Value = -964.1425999999997
$ws.Range("N111").Value = -10674.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2078.907
$ws.Range("I45").Value = 2214.8845
$ws.Range("J45").Value = 1870.9412
$ws.Range("K45").Value = 2214.8845
$ws.Range("L45").Value = 1870.9412
$ws.Range("M45").Value = -1837.8845
$ws.Range("N45").Value = -2624.9412

$ws.Range("H61").Value = 1151.326
$ws.Range("I61").Value = 744.5862
$ws.Range("J61").Value = 1845.1765
$ws.Range("K61").Value = 744.5862
$ws.Range("L61").Value = 1845.1765
$ws.Range("M61").Value = -532.5862
$ws.Range("N61").Value = -2269.1765

$ws.Range("H88").Value = 5090
$ws.Range("I88").Value = 3528.6667
$ws.Range("J88").Value = 5759.143
$ws.Range("K88").Value = 3528.6667
$ws.Range("L88").Value = 5759.143
$ws.Range("M88").Value = -3122.6667
$ws.Range("N88").Value = -6571.143

$ws.Range("H91").Value = 5090
$ws.Range("I91").Value = 3528.6667
$ws.Range("J91").Value = 5759.143
$ws.Range("K91").Value = 3528.6667
$ws.Range("L91").Value = 5759.143
$ws.Range("M91").Value = -2124.6667
$ws.Range("N91").Value = -8567.143

$ws.Range("H122").Value = 23257704
$ws.Range("I122").Value = 29413708
$ws.Range("J122").Value = 1693.5555
$ws.Range("K122").Value = 88241124
$ws.Range("L122").Value = 5080.666499999999
$ws.Range("M122").Value = -88238674
$ws.Range("N122").Value = -9980.666499999999

$ws.Range("H132").Value = 1359.8918
$ws.Range("I132").Value = 1019.82355
$ws.Range("K132").Value = 3059.47065
$ws.Range("M132").Value = -529.4706499999998

$ws.Range("H136").Value = 1151.326
$ws.Range("I136").Value = 744.5862
$ws.Range("J136").Value = 1845.1765
$ws.Range("K136").Value = 2233.7586
$ws.Range("L136").Value = 5535.529500000001
$ws.Range("M136").Value = 316.2413999999999
$ws.Range("N136").Value = -10635.5295

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6712.115
$ws.Range("I134").Value = 953.8823
$ws.Range("J134").Value = 17588.777
$ws.Range("K134").Value = 2861.6469
$ws.Range("L134").Value = 52766.33099999999
$ws.Range("M134").Value = -326.6468999999997
$ws.Range("N134").Value = -57836.33099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71.35294
$ws.Range("J7").Value = 53
$ws.Range("L7").Value = 53
$ws.Range("N7").Value = -279

$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H19").Value = 325.5
$ws.Range("I19").Value = 325.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 325.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -155.5
$ws.Range("N19").ClearContents()

$ws.Range("H23").Value = 5500
$ws.Range("J23").Value = 5500
$ws.Range("L23").Value = 5500
$ws.Range("N23").Value = -5980

$ws.Range("H24").Value = 325.5
$ws.Range("I24").Value = 325.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 325.5
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -155.5
$ws.Range("N24").ClearContents()

$ws.Range("H27").Value = 5500
$ws.Range("J27").Value = 5500
$ws.Range("L27").Value = 5500
$ws.Range("N27").Value = -5884

$ws.Range("H31").Value = 8336515
$ws.Range("I31").Value = 10418188
$ws.Range("J31").Value = 9825
$ws.Range("K31").Value = 10418188
$ws.Range("L31").Value = 9825
$ws.Range("M31").Value = -10417893
$ws.Range("N31").Value = -10415

$ws.Range("H34").Value = 8336515
$ws.Range("I34").Value = 10418188
$ws.Range("J34").Value = 9825
$ws.Range("K34").Value = 10418188
$ws.Range("L34").Value = 9825
$ws.Range("M34").Value = -10417986
$ws.Range("N34").Value = -10229

$ws.Range("H58").Value = 1125.5946
$ws.Range("I58").Value = 1037.9
$ws.Range("J58").Value = 1228.7646
$ws.Range("K58").Value = 1037.9
$ws.Range("L58").Value = 1228.7646
$ws.Range("M58").Value = -834.9000000000001
$ws.Range("N58").Value = -1634.7646

$ws.Range("H122").Value = 1272.579
$ws.Range("I122").Value = 1254.6875
$ws.Range("J122").Value = 1368
$ws.Range("K122").Value = 3764.0625
$ws.Range("L122").Value = 4104
$ws.Range("M122").Value = -1314.0625
$ws.Range("N122").Value = -9004

$ws.Range("H132").Value = 1412.283
$ws.Range("I132").Value = 1307.303
$ws.Range("J132").Value = 1585.5
$ws.Range("K132").Value = 3921.909000000001
$ws.Range("L132").Value = 4756.5
$ws.Range("M132").Value = -1391.909000000001
$ws.Range("N132").Value = -9816.5

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H134").Value = 1186.5278
$ws.Range("I134").Value = 1213.3871
$ws.Range("J134").Value = 1020
$ws.Range("K134").Value = 3640.1613
$ws.Range("L134").Value = 3060
$ws.Range("M134").Value = -1105.1613
$ws.Range("N134").Value = -8130

$ws.Range("H136").Value = 1125.5946
$ws.Range("I136").Value = 1037.9
$ws.Range("J136").Value = 1228.7646
$ws.Range("K136").Value = 3113.7
$ws.Range("L136").Value = 3686.2938
$ws.Range("M136").Value = -563.7000000000003
$ws.Range("N136").Value = -8786.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1115.15
$ws.Range("I5").Value = 1049.8
$ws.Range("J5").Value = 1180.5
$ws.Range("K5").Value = 3149.4
$ws.Range("L5").Value = 3541.5
$ws.Range("M5").Value = -3037.4
$ws.Range("N5").Value = -3765.5

$ws.Range("H113").Value = 954.13635
$ws.Range("I113").Value = 639.4666999999999
$ws.Range("J113").Value = 1046.6863
$ws.Range("K113").Value = 1918.4001
$ws.Range("L113").Value = 3140.0589
$ws.Range("M113").Value = 251.5999000000002
$ws.Range("N113").Value = -7480.0589

$ws.Range("H135").Value = 1115.15
$ws.Range("I135").Value = 1049.8
$ws.Range("J135").Value = 1180.5
$ws.Range("K135").Value = 9448.199999999999
$ws.Range("L135").Value = 10624.5
$ws.Range("M135").Value = -6913.199999999999
$ws.Range("N135").Value = -15694.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2766.4827
$ws.Range("I7").Value = 2120.5
$ws.Range("J7").Value = 3561.5386
$ws.Range("K7").Value = 2120.5
$ws.Range("L7").Value = 3561.5386
$ws.Range("M7").Value = -2008.5
$ws.Range("N7").Value = -3785.5386

$ws.Range("H17").Value = 150495
$ws.Range("I17").Value = 150495
$ws.Range("K17").Value = 150495
$ws.Range("M17").Value = -150325

$ws.Range("H18").Value = 1552.5
$ws.Range("I18").Value = 1552.5
$ws.Range("K18").Value = 1552.5
$ws.Range("M18").Value = -1380.5

$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -830

$ws.Range("H20").Value = 3000
$ws.Range("J20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3452

$ws.Range("H23").Value = 40000000
$ws.Range("I23").Value = 40000000
$ws.Range("K23").Value = 40000000
$ws.Range("M23").Value = -39999770

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H40").Value = 1913.0769
$ws.Range("I40").Value = 1704.2858
$ws.Range("J40").Value = 2156.6667
$ws.Range("K40").Value = 1704.2858
$ws.Range("L40").Value = 2156.6667
$ws.Range("M40").Value = -1568.2858
$ws.Range("N40").Value = -2428.6667

$ws.Range("H122").Value = 2614.4285
$ws.Range("I122").Value = 2533.7778
$ws.Range("J122").Value = 2759.6
$ws.Range("K122").Value = 7601.3334
$ws.Range("L122").Value = 8278.799999999999
$ws.Range("M122").Value = -5151.3334
$ws.Range("N122").Value = -13178.8

$ws.Range("H126").Value = 2766.4827
$ws.Range("I126").Value = 2120.5
$ws.Range("J126").Value = 3561.5386
$ws.Range("K126").Value = 6361.5
$ws.Range("L126").Value = 10684.6158
$ws.Range("M126").Value = -3891.5
$ws.Range("N126").Value = -15624.6158

$ws.Range("H136").Value = 2342.9333
$ws.Range("I136").Value = 1685.591
$ws.Range("J136").Value = 4150.625
$ws.Range("K136").Value = 5056.772999999999
$ws.Range("L136").Value = 12451.875
$ws.Range("M136").Value = -2506.772999999999
$ws.Range("N136").Value = -17551.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 70007
$ws.Range("J12").Value = 70007
$ws.Range("L12").Value = 70007
$ws.Range("N12").Value = -70291

$ws.Range("H13").Value = 14999.333
$ws.Range("I13").Value = 12500
$ws.Range("J13").Value = 19998
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 19998
$ws.Range("M13").Value = -12360
$ws.Range("N13").Value = -20278

$ws.Range("H107").Value = 496.10526
$ws.Range("I107").Value = 502.53333
$ws.Range("J107").Value = 472
$ws.Range("K107").Value = 1507.59999
$ws.Range("L107").Value = 1416
$ws.Range("M107").Value = 412.4000100000001
$ws.Range("N107").Value = -5256

$ws.Range("H136").Value = 580.1111
$ws.Range("I136").Value = 423.75
$ws.Range("J136").Value = 837.64703
$ws.Range("K136").Value = 1271.25
$ws.Range("L136").Value = 2512.94109
$ws.Range("M136").Value = 1278.75
$ws.Range("N136").Value = -7612.94109
